$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Snapshot the original data (before shifting anything) ---
# Headers (row 1): old B..F -> text for CUT_SPEED, PIERCE_TIME, WEIGHT,
# MATERIAL_COST, COST_SQUARE_FOOT
$hdrCutSpeed   = $ws.Cells.Item(1, 2).Value2
$hdrPierceTime = $ws.Cells.Item(1, 3).Value2
$hdrWeight     = $ws.Cells.Item(1, 4).Value2
$hdrMatCost    = $ws.Cells.Item(1, 5).Value2
$hdrCostSqFt   = $ws.Cells.Item(1, 6).Value2

# Data rows 2..7, old columns B..E (numeric data; old F was the formula,
# rebuilt fresh below instead of copied)
$cutSpeed   = @{}
$pierceTime = @{}
$weight     = @{}
$matCost    = @{}
for ($r = 2; $r -le 7; $r++) {
    $cutSpeed[$r]   = $ws.Cells.Item($r, 2).Value2
    $pierceTime[$r] = $ws.Cells.Item($r, 3).Value2
    $weight[$r]     = $ws.Cells.Item($r, 4).Value2
    $matCost[$r]    = $ws.Cells.Item($r, 5).Value2
}

# New thickness values (numeric) for column A, rows 2..7 (replacing the old
# "0.035 [20GA]"-style text labels)
$thickness = @{ 2 = 0.035; 3 = 0.047; 4 = 0.06; 5 = 0.075; 6 = 0.12; 7 = 0.187 }

# --- Clear the whole used range and rebuild it in the new layout ---
$ws.Range("A1:F7").ClearContents()

# Row 1 headers: A1=STEEL, C1=CUT_SPEED, D1=PIERCE_TIME, E1=WEIGHT,
# F1=MATERIAL_COST, G1=COST_SQUARE_FOOT (B1 stays blank - new grade column)
$ws.Cells.Item(1, 1).Value = "STEEL"
$ws.Cells.Item(1, 3).Value = $hdrCutSpeed
$ws.Cells.Item(1, 4).Value = $hdrPierceTime
$ws.Cells.Item(1, 5).Value = $hdrWeight
$ws.Cells.Item(1, 6).Value = $hdrMatCost
$ws.Cells.Item(1, 7).Value = $hdrCostSqFt

# Data rows 2..7
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = $thickness[$r]
    $ws.Cells.Item($r, 2).Value = "s"
    $ws.Cells.Item($r, 3).Value = $cutSpeed[$r]
    $ws.Cells.Item($r, 4).Value = $pierceTime[$r]
    $ws.Cells.Item($r, 5).Value = $weight[$r]
    $ws.Cells.Item($r, 6).Value = $matCost[$r]
}

# Cost/sqft formula, rebuilt as one shared formula across G2:G7
$ws.Range("G2:G7").Formula = "=E2*F2"

# --- Number formats ---
$ws.Range("A2:A7").NumberFormat = "0.000"
$ws.Range("B2:B7").NumberFormat = "0.000"
$ws.Range("D2:D7").NumberFormat = "0.0"
$ws.Range("E2:F7").NumberFormat = "0.00"
$ws.Range("G2:G7").NumberFormat = "0.000"

# --- Extra blank rows 8..12 (columns A/B only), same number format ---
$ws.Range("A8:B12").NumberFormat = "0.000"

# --- Column width for the new grade column B ---
$ws.Columns.Item(2).ColumnWidth = 7.85546875

# --- Selection to match the saved view ---
$ws.Range("A2:G7").Select()

# --- Minimize the window (matches workbookView minimized="1") ---
$excel.ActiveWindow.WindowState = -4140  # xlMinimized
